# Version 3: Agregar texto a 2 archivos.
#
# Recreates, via the Excel object model, the edit that:
#   - adds a shared string "AGREGANDO CAMBIOS AL REPOSITORIO"
#   - adds a bold, 14pt, centered title style
#   - writes that title into B2:I2 on Hoja1, merges the range,
#     sets the row height, and leaves it selected
#   - sets the page orientation to portrait

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the title's look (bold, 14pt, centered) on a throw-away cell far off
# the used area. Doing the formatting here - instead of directly on B2:I2 -
# means only a single "paste formats" operation touches the title row, so
# we don't leave a trail of partially-formatted intermediate styles on the
# cells that actually matter.
$scratch = $ws.Range("Z100")
$scratch.Font.Bold = $true
$scratch.Font.Size = 14
$scratch.HorizontalAlignment = -4108   # xlCenter

$titleRange = $ws.Range("B2:I2")

# Apply that formatting (and only the formatting) to the whole title range.
$scratch.Copy()
$titleRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$scratch.Clear() | Out-Null

# Title text goes in the top-left cell of the range, then merge the row.
$titleRange.Value = "AGREGANDO CAMBIOS AL REPOSITORIO"
$titleRange.Merge() | Out-Null

# Taller title row, matching the source workbook.
$ws.Rows.Item(2).RowHeight = 18.75

# Page orientation, as set in the source workbook.
$ws.PageSetup.Orientation = 1   # xlPortrait

# Leave the merged title selected/active, like in the source workbook.
$titleRange.Select() | Out-Null
